# Apply "stock splits also for ISIN" update:
# add 6 new rows (r=8..13) to Sheet1 holding the same split history as
# NVDA's ticker rows (2..7), but keyed by the NVDA ISIN "US67066G1040"
# instead of the symbol, using a slightly different (explicit) font color.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy existing formatting for each column down into the new rows so the
# new cells pick up the same number format / font as rows 2:7.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A8:A13").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B8:B13").PasteSpecial(-4122) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("C8:C13").PasteSpecial(-4122) | Out-Null

# New split history rows, identical values to the NVDA symbol rows above,
# but identified by the NVDA ISIN.
$rows = @(
    @{r = 8;  date = 45450; shares = 10},
    @{r = 9;  date = 44397; shares = 4},
    @{r = 10; date = 39336; shares = 1.5},
    @{r = 11; date = 38814; shares = 2},
    @{r = 12; date = 37146; shares = 2},
    @{r = 13; date = 36704; shares = 2}
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.r, 1).Value2 = $row.date
    $ws.Cells.Item($row.r, 2).Value2 = "US67066G1040"
    $ws.Cells.Item($row.r, 3).Value2 = $row.shares
}

# Give the ISIN column its own (slightly different) font color, which is
# what introduces the extra font / cellXf entry in styles.xml.
$ws.Range("B8:B13").Font.Name = "Calibri"
$ws.Range("B8:B13").Font.Size = 12
$ws.Range("B8:B13").Font.Color = 2630431

# Reset the selection back to the top-left cell.
$ws.Range("A1").Select() | Out-Null
